$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range('D2').Value = '57.320.02'
$ws.Range("E2").NumberFormat = "@"
$ws.Range('E2').Value = '  +1.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range('D3').Value = '3.024.60'
$ws.Range("E3").NumberFormat = "@"
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '519.45'
$ws.Range("E5").NumberFormat = "@"
$ws.Range('E5').Value = '  +5.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '141.02'
$ws.Range("E6").NumberFormat = "@"
$ws.Range('E6').Value = '  +5.25%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range('E8').Value = '  +3.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '7.60'
$ws.Range("E9").NumberFormat = "@"
$ws.Range('E9').Value = '  +5.43%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range('E10').Value = '  +6.22%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range('E12').Value = '  +2.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '3.540.85'
$ws.Range("E13").NumberFormat = "@"
$ws.Range('E13').Value = '  +1.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '26.15'
$ws.Range("E14").NumberFormat = "@"
$ws.Range('E14').Value = '  +5.11%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range('E15').Value = '  +12.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '57.309.58'
$ws.Range("E16").NumberFormat = "@"
$ws.Range('E16').Value = '  +1.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '3.024.16'
$ws.Range("E17").NumberFormat = "@"
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '6.04'
$ws.Range("E18").NumberFormat = "@"
$ws.Range('E18').Value = '  +3.41%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range('E19').Value = '  +3.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '7.98'
$ws.Range("E20").NumberFormat = "@"
$ws.Range('E20').Value = '  +3.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '332.60'
$ws.Range("E21").NumberFormat = "@"
$ws.Range('E21').Value = '  +3.18%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '0.489'
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '64.05'
$ws.Range("E24").NumberFormat = "@"
$ws.Range('E24').Value = '  +5.05%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range('E25').Value = '  +5.96%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '0.0₃0928'
$ws.Range("E27").NumberFormat = "@"
$ws.Range('E27').Value = '  +6.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '6.75'
$ws.Range("E28").NumberFormat = "@"
$ws.Range('E28').Value = '  +3.46%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range('E29').Value = '  +7.17%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range("C30").NumberFormat = "@"
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '1.23'
$ws.Range("E30").NumberFormat = "@"
$ws.Range('E30').Value = '  +5.17%  '
$ws.Range("B31").NumberFormat = "@"
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range("C31").NumberFormat = "@"
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '1.83'
$ws.Range("E31").NumberFormat = "@"
$ws.Range('E31').Value = '  +6.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '20.84'
$ws.Range("E32").NumberFormat = "@"
$ws.Range('E32').Value = '  +5.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '157.84'
$ws.Range("E33").NumberFormat = "@"
$ws.Range('E33').Value = '  +4.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '4.65'
$ws.Range("E34").NumberFormat = "@"
$ws.Range('E34').Value = '  +4.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '5.79'
$ws.Range("E35").NumberFormat = "@"
$ws.Range('E35').Value = '  +2.89%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range('E36').Value = '  +1.74%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range('B37').Value = 'Hedera'
$ws.Range("C37").NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '0.0682'
$ws.Range("E37").NumberFormat = "@"
$ws.Range('E37').Value = '  +3.28%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range("C38").NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '24.37'
$ws.Range("E38").NumberFormat = "@"
$ws.Range('E38').Value = '  +3.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '3.055.52'
$ws.Range("E39").NumberFormat = "@"
$ws.Range('E39').Value = '  +0.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '37.39'
$ws.Range("E40").NumberFormat = "@"
$ws.Range('E40').Value = '  +1.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range("E41").NumberFormat = "@"
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '2.300.97'
$ws.Range("E42").NumberFormat = "@"
$ws.Range('E42').Value = '  +6.07%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range('E43').Value = '  +1.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '3.75'
$ws.Range("E44").NumberFormat = "@"
$ws.Range('E44').Value = '  +5.65%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range('E45').Value = '  +2.20%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range('E47').Value = '  +10.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '0.0243'
$ws.Range("E48").NumberFormat = "@"
$ws.Range('E48').Value = '  +2.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '5.92'
$ws.Range("E49").NumberFormat = "@"
$ws.Range('E49').Value = '  +6.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '19.53'
$ws.Range("E50").NumberFormat = "@"
$ws.Range('E50').Value = '  +1.11%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range('E51').Value = '  +4.25%  '
